$p = $ppt.ActivePresentation

# --- Slide 5 ("Why did you choose to learn web development?") ---
# "Entrepreneurial Ventures: " paragraph: the run containing just a
# trailing space and the run containing just "W" (both plain,
# non-bold "en-US" runs) get merged into a single " W" run, i.e. the
# text is untouched but the run split collapses.
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2).TextFrame.TextRange

$found = $body5.Find("Ventures: W", 0)
if ($found -ne $null) {
    $spacePos = $found.Start + $found.Length - 2
    $pair = $body5.Characters($spacePos, 2)
    $pair.Text = " W"
}

# --- Slide 2 ("How does the web work?") ---
# The "Content Placeholder 2" body box grew taller (its bottom moved
# down), which is what drives PowerPoint's autofit to recompute.
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2)
$body2.Height = 313.9336220472441
